$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly added UPE numbers for rows 2 and 3 (column A)
$ws.Range("A2").Value = 188292
$ws.Range("A3").Value = 105424
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Style = "Normal"

# Update the active cell selection to A3, matching the post-edit state
$ws.Range("A3").Select()
